# Applies the commit "added function to get list and default values":
#   1. For a set of "display" impression rows, rename the impression id in
#      column A from "...display" to "...video" and bump the column C
#      count value from 1 to 2 (these rows become video impressions).
#   2. For a set of rows whose VIDEO_STARTED (I) / VIDEO_MRC_VIEWED (J)
#      columns were stored as plain numbers (0/1), convert those cells to
#      proper boolean values (TRUE/FALSE) while preserving the same
#      underlying 0/1 semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column-A impression id must change from "display" to "video",
# with the corresponding column C value changing from 1 to 2.
$renameRows = @(39, 49, 50, 51, 52, 53, 54, 55, 56, 62, 70, 80, 107, 139, 141, 147)

foreach ($r in $renameRows) {
    $aCell = $ws.Cells.Item($r, 1)
    $oldName = $aCell.Value2
    if ($oldName -like "*display") {
        $newName = $oldName -replace "display", "video"
        $aCell.Value = $newName
    }
    $ws.Cells.Item($r, 3).Value = 2
}

# Rows whose I (VIDEO_STARTED) and J (VIDEO_MRC_VIEWED) columns are stored
# as numbers and must be converted to booleans, keeping the same value.
$boolRows = @(38, 40, 41, 42, 43, 61, 63, 64, 65, 66, 67, 68, 69, 71, 72, 73, 74, 75, 76, 77, 78, 79, 81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 96, 106, 108, 109, 110, 111, 132, 133, 134, 135, 136, 137, 138, 140, 142, 143, 144, 145, 146)

foreach ($r in $boolRows) {
    $iCell = $ws.Cells.Item($r, 9)
    $jCell = $ws.Cells.Item($r, 10)

    $iVal = $iCell.Value2
    $jVal = $jCell.Value2

    if ($iVal -eq 1) {
        $iCell.Value = $true
    } else {
        $iCell.Value = $false
    }

    if ($jVal -eq 1) {
        $jCell.Value = $true
    } else {
        $jCell.Value = $false
    }
}
